# Update Yearn partnership metrics on the "Strategic Partnership - Yearn Finance" slide.
# TVL $804M(.21M) -> $965M(.92M), vault counts 18 -> 33.

$p = $ppt.ActivePresentation

# Locate the slide/shape that holds the Yearn partnership copy (searches all slides so the
# script is resilient to slide re-ordering).
$targetShape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -like "*804M+ TVL*") {
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Simple whole-run replacements: old/new text is plain ASCII, so a straight
# IndexOf/Characters round-trip is safe.
$replacements = @(
    @(" (`$804M+ TVL) and ", " (`$965M+ TVL) and "),
    @(" Partnership with #1 DeFi yield aggregator (`$804M+ TVL) - ", " Partnership with #1 DeFi yield aggregator (`$965M+ TVL) - "),
    @(" 18 battle-tested vaults to template and compose", " 33 battle-tested vaults to template and compose"),
    @("18 Yearn Vaults", "33 Yearn Vaults"),
    @("`$804.21M TVL", "`$965.92M TVL")
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]

    $current = $tr.Text
    $idx = $current.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "Could not find expected run text: $oldText"
    }

    $sub = $tr.Characters($idx + 1, $oldText.Length)
    $sub.Text = $newText
}

# The "For Users:" bullet run contains a non-breaking space (U+00A0, "vs. months") that the
# COM text getter mangles when read back, so locate it via an ASCII-only prefix and use the
# known run length (230 UTF-16 code units) instead of round-tripping the nbsp through IndexOf.
$nbspPrefix = " - Create custom Yearn strategies without coding"
$nbspRunLength = 230
$current = $tr.Text
$idx = $current.IndexOf($nbspPrefix)
if ($idx -lt 0) {
    throw "Could not find expected run text: $nbspPrefix"
}
$sub = $tr.Characters($idx + 1, $nbspRunLength)
$nbsp = [char]0x00A0
$sub.Text = " - Create custom Yearn strategies without coding - Charge your own management fees (you become the vault operator) - Access `$965M+ Yearn ecosystem TVL and proven strategies - Deploy in minutes vs." + $nbsp + "months of traditional development"
